$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column letters that receive new values (A, D, L, M are unchanged).
$cols = @("B","C","E","F","G","H","I","J","K","N","O")

# New values for rows 2..25, one array per row, in the same order as $cols.
$newValues = @(
        @(8.330358067558683,5.245086243497296,24.71350566460032,39.16134544342578,23.1509785461277,13.01325835140012,18.38388474352791,7.645351252191472,8.066406362232499,17.2477576635746,19.03689494036477),
        @(8.000109851581756,5.084424757290152,24.39100030984826,38.97303980010369,23.2559722498918,13.05969800669259,18.47364334783031,7.664930123680536,7.773221910830254,17.29076378730346,19.11854880921704),
        @(7.790872487564901,4.982386937699148,24.1969539323007,38.86781928162924,23.32807108760148,13.09005124820734,18.53198862871566,7.677648911843264,7.622735934146376,17.31867845268268,19.17246895364721),
        @(7.704104046606272,4.939994029706461,24.11897057250206,38.8275897324293,23.35936121431279,13.10288331035658,18.55657855354491,7.683007675996249,7.564303078295577,17.33043412981202,19.1953924881433),
        @(7.689609337325827,4.932906969882729,24.10609013125547,38.82107051855483,23.36467193549046,13.1050420330416,18.56071086297759,7.683908122170872,7.554566067345839,17.33240914057174,19.19925631200786),
        @(7.789708208794845,4.981818441546268,24.19589767826425,38.86726596633494,23.32848535998968,13.09022243119219,18.53231696118645,7.677720469829391,7.621950242143622,17.31883545322093,19.1727742606222),
        @(8.21789929421136,5.190417424010424,24.6015375974401,39.0942805011437,23.18559040275493,13.02888928504723,18.41416292223541,7.651957577598512,7.966864340820921,17.26227354351739,19.06426317251054),
        @(9.001604659990042,5.570959895676215,25.42435025377697,39.62031214474187,22.96636444455567,12.92318896195629,18.20808225926502,7.606950264924009,8.655559722888032,17.1632887382745,18.88154361623148),
        @(9.537897479147201,5.831251201834129,26.03989778419768,40.0535961211289,22.8430519633223,12.8543895274177,18.07224456224639,7.577218198111928,9.121674202742664,17.09778611079511,18.7656908724834),
        @(9.772394112815807,5.945139326873177,26.32112373693212,40.26027012062145,22.79526684771253,12.82500920199449,18.01381974862309,7.564410849514546,9.324547642549595,17.06954423615808,18.71699329810553),
        @(9.859769220186093,5.987591903372323,26.4276922523951,40.33985309294153,22.77837576537417,12.81415894918124,17.99217961107486,7.559663859255843,9.400017195412929,17.05907265377757,18.69912985303071),
        @(9.84101564711715,5.978479334392341,26.40473907370212,40.32265569224739,22.78195986945865,12.81648349815767,17.9968186770583,7.560681638858842,9.383824231144706,17.06131798914071,18.70295136818738),
        @(9.779611374585762,5.948645571626353,26.32989026624158,40.26679133643045,22.79385302572552,12.82411102608628,18.01202970155853,7.564018252056901,9.330783898820446,17.0686782673719,18.71551208407003),
        @(9.74181237678035,5.930283028140146,26.28404997009488,40.23274299424335,22.80129499921007,12.82881897159646,18.02140992537282,7.566075410676744,9.298117812665216,17.07321567159284,18.72328110582988),
        @(9.522376546980482,5.823715241198824,26.02153479357109,40.0402773328485,22.84634284529466,12.85634815531642,18.07613051559077,7.578069603830566,9.108228305987076,17.09966302697244,18.76895407514462),
        @(9.385287985108496,5.757163436697043,25.86072867723988,39.92461807210636,22.87611332395992,12.8737272251796,18.11056233060012,7.585611270360613,8.989362463449044,17.11628555297228,18.79799960711925),
        @(9.305550441272349,5.718460456542417,25.76835702819466,39.85899833992917,22.89401820896851,12.88390364382688,18.13068366842419,7.590016632561714,8.920133204146641,17.12599285095548,18.81508265632377),
        @(9.278402248402005,5.70528423783937,25.73710521642975,39.83693763001735,22.90021446789369,12.88738020118694,18.13755087457979,7.591519833110805,8.896546722949052,17.12930475131764,18.82093136027653),
        @(9.39997366506293,5.764292073247895,25.8778350950257,39.9368369932049,22.87286323614008,12.87185852048062,18.10686419206149,7.584801453877003,9.002105301849399,17.11450090361754,18.79486865006567),
        @(9.797686390225298,5.957426954463791,26.35187394150016,40.28316468881727,22.79032696822371,12.8218631624738,18.00754871959077,7.563035418742754,9.346400147783415,17.0665103284426,18.71180701883657),
        @(10.04929344547667,6.079711535366966,26.66207446579965,40.51717599716936,22.74340707537039,12.79079368990525,17.94546174663244,7.549409525008293,9.563509950341114,17.03644523334176,18.66088689227928),
        @(9.915785812227396,6.014813844895045,26.4965114018025,40.39159757026731,22.76780367770814,12.80722922333421,17.97834066059008,7.556627187514771,9.448368556285439,17.05237285857551,18.68775550218331),
        @(9.39333714081571,5.761070589860715,25.87010103313656,39.93131009121625,22.87433014246624,12.87270278571144,18.10853510490189,7.585167355215559,8.996347039281574,17.11530727384182,18.79628295850442),
        @(8.796185898812478,5.471276535536501,25.19941174175739,39.46959244178716,23.01908109726532,12.9502262232632,18.26109499813772,7.618538415504567,8.476055071961024,17.18879479629128,18.92774879806168)
)

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $newValues[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range($cols[$j] + $rowNum).Value = $rowVals[$j]
    }
}
